$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: header "time_taken" in F1, formatted like the other header cells (E1)
$ws.Cells.Item(1, 6).Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Timestamp values for rows 2-31 (time_taken metadata), unformatted like the rest of the data cells
$timestamps = @(
    "2021-10-05 13:41:11.412492",
    "2021-10-05 13:41:11.412504",
    "2021-10-05 13:41:11.412508",
    "2021-10-05 13:41:11.412511",
    "2021-10-05 13:41:11.412514",
    "2021-10-05 13:41:11.412518",
    "2021-10-05 13:41:11.412520",
    "2021-10-05 13:41:11.412523",
    "2021-10-05 13:41:11.412527",
    "2021-10-05 13:41:11.412598",
    "2021-10-05 13:41:11.412603",
    "2021-10-05 13:41:11.412606",
    "2021-10-05 13:41:11.412609",
    "2021-10-05 13:41:11.412611",
    "2021-10-05 13:41:11.412636",
    "2021-10-05 13:41:11.412643",
    "2021-10-05 13:41:11.412646",
    "2021-10-05 13:41:11.412649",
    "2021-10-05 13:41:11.412652",
    "2021-10-05 13:41:11.412654",
    "2021-10-05 13:41:11.412657",
    "2021-10-05 13:41:11.412659",
    "2021-10-05 13:41:11.412662",
    "2021-10-05 13:41:11.412664",
    "2021-10-05 13:41:11.412667",
    "2021-10-05 13:41:11.412670",
    "2021-10-05 13:41:11.412672",
    "2021-10-05 13:41:11.412675",
    "2021-10-05 13:41:11.412677",
    "2021-10-05 13:41:11.412680"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
